# Replace the (outdated, Perseus/2018) "Kampagnendaten" campaign-date
# sentence with the translated Bootes dates in every paragraph that still
# carries the old text. The new sentence is emitted as a single, plain
# (un-formatted) run - matching how Word collapses a fully
# select-and-retype edit into one <w:r><w:t>...</w:t></w:r> with no
# <w:rPr> - so each paragraph's run list is first cleared down to zero
# runs and the replacement text is written in a second pass.

$d = $word.ActiveDocument
$newText = "Kampagnendaten Bootes: 14. bis 23. Mai, 13. bis 22. Juni, 12.-21. Juli"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "Kampagnendaten*") {
        # Range covering just the paragraph's text, excluding the
        # trailing paragraph-mark character (End is one past it).
        $start = $p.Range.Start
        $end = $p.Range.End - 1

        # Pass 1: wipe every run in the paragraph (leaves zero runs).
        $d.Range($start, $end).Text = ""

        # Pass 2: write the new sentence into the now run-less
        # paragraph - produces a single fresh run with no inherited
        # rPr/rsid/xml:space baggage.
        $d.Paragraphs.Item($i).Range.Text = $newText
    }
}
